# Update cumulative season target-depth stats on both the OFF and DEF
# sheets to include the simulated Wild Card round game.

$wb = $excel.ActiveWorkbook

$offSheet = $wb.Worksheets.Item("OFF")
$offSheet.Range("B3").Value = 231
$offSheet.Range("C3").Value = 162
$offSheet.Range("D3").Value = 63
$offSheet.Range("E3").Value = 37
$offSheet.Range("F3").Value = 4
$offSheet.Range("G3").Value = 4

$defSheet = $wb.Worksheets.Item("DEF")
$defSheet.Range("B3").Value = 231
$defSheet.Range("C3").Value = 171
$defSheet.Range("D3").Value = 52
$defSheet.Range("E3").Value = 25
$defSheet.Range("F3").Value = 2
$defSheet.Range("G3").Value = 3
